# WBS Who - apply "World barriers in, bug fixed with dungeons, sound effects in"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mark Medieval World / Dungeon branch as Done ---
$ws.Range("B37").Value = "Done"   # Medieval World
$ws.Range("B38").Value = "Done"   # Medieval World - Design
$ws.Range("B39").Value = "Done"   # Medieval World - Implement
$ws.Range("B40").Value = "Done"   # Dungeon

# --- Ai now also worked on by Joe ---
$ws.Range("C43").Value = "Alex/ Joe"   # Ai - Who

# --- Insert a new sub-task row for Ai pathfinding, right below "Ai" (row 43) ---
$ws.Rows("44:44").Insert()
$ws.Range("A44").Value = "Ai - PathFinding"
$ws.Range("B44").Value = "Started"
$ws.Range("C44").Value = "Joe"

# --- Everything below shifted down by one row; apply remaining content edits ---
$ws.Range("C53").Value = "Matt/Alex"   # NPC - Who
$ws.Range("B54").Value = "Done"        # Character Saves/Loads - Progress

$ws.Range("C58").Value = "John"        # Chips sold for Exp - Who
$ws.Range("C59").Value = "John"        # Pickup items - Who

$ws.Range("B61").Value = "Done"         # Potions - Progress
$ws.Range("C61").Value = "Alex"         # Potions - Who
$ws.Range("G61").Value = "Behind- Done" # Potions - Goal Met

$ws.Range("B63").Value = "Done"  # Magic Effects - Progress
$ws.Range("C63").Value = "John"  # Magic Effects - Who
$ws.Range("G63").Value = "Done"  # Magic Effects - Goal Met

$ws.Range("C66").Value = "Alex"  # World Barriers - Who

# --- Update the saved view state to match the author's last position ---
$ws.Application.ActiveWindow.ScrollRow = 31
$ws.Range("B39").Select()
